$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$arr = New-Object 'object[,]' 25,16
$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 10.070632
$arr[0,3] = 30.211896
$arr[0,4] = 0.07634150103324112
$arr[0,5] = 0.08026042296304617
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 10.68421466666667
$arr[0,9] = 32.052644
$arr[0,10] = 0.04705285980693976
$arr[0,11] = 0.04892736897547583
$arr[0,12] = 107.5967941170027
$arr[0,13] = 968.3711470530241
$arr[0,14] = 0.003592085945568441
$arr[0,15] = 0.003926931328440713
$arr[1,0] = 3
$arr[1,1] = 1
$arr[1,2] = 10.070632
$arr[1,3] = 30.211896
$arr[1,4] = 0.07634150103324112
$arr[1,5] = 0.08026042296304617
$arr[1,6] = 3
$arr[1,7] = 1
$arr[1,8] = 70.36235166666667
$arr[1,9] = 211.087055
$arr[1,10] = 0.3098730203341347
$arr[1,11] = 0.3222178559101571
$arr[1,12] = 708.5933502895866
$arr[1,13] = 6377.34015260628
$arr[1,14] = 0.02365617150201189
$arr[1,15] = 0.02586134140159507
$arr[2,0] = 3
$arr[2,1] = 1
$arr[2,2] = 10.070632
$arr[2,3] = 30.211896
$arr[2,4] = 0.07634150103324112
$arr[2,5] = 0.08026042296304617
$arr[2,6] = 3
$arr[2,7] = 1
$arr[2,8] = 59.09107466666666
$arr[2,9] = 177.273224
$arr[2,10] = 0.2602347611759026
$arr[2,11] = 0.2706020894912812
$arr[2,12] = 595.0844674525226
$arr[2,13] = 5355.760207072703
$arr[2,14] = 0.01986671228919543
$arr[2,15] = 0.0217186381572543
$arr[3,0] = 3
$arr[3,1] = 1
$arr[3,2] = 10.070632
$arr[3,3] = 30.211896
$arr[3,4] = 0.07634150103324112
$arr[3,5] = 0.08026042296304617
$arr[3,6] = 3
$arr[3,7] = 1
$arr[3,8] = 60.83231733333333
$arr[3,9] = 182.496952
$arr[3,10] = 0.2679031251727568
$arr[3,11] = 0.2785759485989269
$arr[3,12] = 612.6198815712213
$arr[3,13] = 5513.578934140992
$arr[3,14] = 0.02045212670718454
$arr[3,15] = 0.02235862346188168
$arr[4,0] = 3
$arr[4,1] = 1
$arr[4,2] = 10.070632
$arr[4,3] = 30.211896
$arr[4,4] = 0.07634150103324112
$arr[4,5] = 0.08026042296304617
$arr[4,6] = 2
$arr[4,7] = 1
$arr[4,8] = 26.0983795
$arr[4,9] = 52.196759
$arr[4,10] = 0.1149362335102661
$arr[4,11] = 0.07967673702415903
$arr[4,12] = 262.827175740844
$arr[4,13] = 1576.963054445064
$arr[4,14] = 0.008774404589280827
$arr[4,15] = 0.006394888613874405
$arr[5,0] = 3
$arr[5,1] = 1
$arr[5,2] = 42.14988333333334
$arr[5,3] = 126.44965
$arr[5,4] = 0.3195216905992255
$arr[5,5] = 0.3359240476840365
$arr[5,6] = 3
$arr[5,7] = 1
$arr[5,8] = 10.68421466666667
$arr[5,9] = 32.052644
$arr[5,10] = 0.04705285980693976
$arr[5,11] = 0.04892736897547583
$arr[5,12] = 450.3384017082889
$arr[5,13] = 4053.0456153746
$arr[5,14] = 0.01503440931304174
$arr[5,15] = 0.01643587982877219
$arr[6,0] = 3
$arr[6,1] = 1
$arr[6,2] = 42.14988333333334
$arr[6,3] = 126.44965
$arr[6,4] = 0.3195216905992255
$arr[6,5] = 0.3359240476840365
$arr[6,6] = 3
$arr[6,7] = 1
$arr[6,8] = 70.36235166666667
$arr[6,9] = 211.087055
$arr[6,10] = 0.3098730203341347
$arr[6,11] = 0.3222178559101571
$arr[6,12] = 2965.764913808972
$arr[6,13] = 26691.88422428075
$arr[6,14] = 0.09901115132825092
$arr[6,15] = 0.1082407263934116
$arr[7,0] = 3
$arr[7,1] = 1
$arr[7,2] = 42.14988333333334
$arr[7,3] = 126.44965
$arr[7,4] = 0.3195216905992255
$arr[7,5] = 0.3359240476840365
$arr[7,6] = 3
$arr[7,7] = 1
$arr[7,8] = 59.09107466666666
$arr[7,9] = 177.273224
$arr[7,10] = 0.2602347611759026
$arr[7,11] = 0.2706020894912812
$arr[7,12] = 2490.681903241289
$arr[7,13] = 22416.1371291716
$arr[7,14] = 0.0831506508436101
$arr[7,15] = 0.09090174921366907
$arr[8,0] = 3
$arr[8,1] = 1
$arr[8,2] = 42.14988333333334
$arr[8,3] = 126.44965
$arr[8,4] = 0.3195216905992255
$arr[8,5] = 0.3359240476840365
$arr[8,6] = 3
$arr[8,7] = 1
$arr[8,8] = 60.83231733333333
$arr[8,9] = 182.496952
$arr[8,10] = 0.2679031251727568
$arr[8,11] = 0.2785759485989269
$arr[8,12] = 2564.075078496311
$arr[8,13] = 23076.6757064668
$arr[8,14] = 0.08560085947201519
$arr[8,15] = 0.09358036024077161
$arr[9,0] = 3
$arr[9,1] = 1
$arr[9,2] = 42.14988333333334
$arr[9,3] = 126.44965
$arr[9,4] = 0.3195216905992255
$arr[9,5] = 0.3359240476840365
$arr[9,6] = 2
$arr[9,7] = 1
$arr[9,8] = 26.0983795
$arr[9,9] = 52.196759
$arr[9,10] = 0.1149362335102661
$arr[9,11] = 0.07967673702415903
$arr[9,12] = 1100.043651114058
$arr[9,13] = 6600.26190668435
$arr[9,14] = 0.0367246196423076
$arr[9,15] = 0.02676533200741204
$arr[10,0] = 3
$arr[10,1] = 1
$arr[10,2] = 26.941003
$arr[10,3] = 80.823009
$arr[10,4] = 0.2042291495073052
$arr[10,5] = 0.2147130682392819
$arr[10,6] = 3
$arr[10,7] = 1
$arr[10,8] = 10.68421466666667
$arr[10,9] = 32.052644
$arr[10,10] = 0.04705285980693976
$arr[10,11] = 0.04892736897547583
$arr[10,12] = 287.8434593873107
$arr[10,13] = 2590.591134485796
$arr[10,14] = 0.009609565540257773
$arr[10,15] = 0.01050534551359987
$arr[11,0] = 3
$arr[11,1] = 1
$arr[11,2] = 26.941003
$arr[11,3] = 80.823009
$arr[11,4] = 0.2042291495073052
$arr[11,5] = 0.2147130682392819
$arr[11,6] = 3
$arr[11,7] = 1
$arr[11,8] = 70.36235166666667
$arr[11,9] = 211.087055
$arr[11,10] = 0.3098730203341347
$arr[11,11] = 0.3222178559101571
$arr[11,12] = 1895.632327338722
$arr[11,13] = 17060.69094604849
$arr[11,14] = 0.06328510339810023
$arr[11,15] = 0.06918438448395264
$arr[12,0] = 3
$arr[12,1] = 1
$arr[12,2] = 26.941003
$arr[12,3] = 80.823009
$arr[12,4] = 0.2042291495073052
$arr[12,5] = 0.2147130682392819
$arr[12,6] = 3
$arr[12,7] = 1
$arr[12,8] = 59.09107466666666
$arr[12,9] = 177.273224
$arr[12,10] = 0.2602347611759026
$arr[12,11] = 0.2706020894912812
$arr[12,12] = 1591.97281986789
$arr[12,13] = 14327.75537881102
$arr[12,14] = 0.05314752394719129
$arr[12,15] = 0.05810180490663373
$arr[13,0] = 3
$arr[13,1] = 1
$arr[13,2] = 26.941003
$arr[13,3] = 80.823009
$arr[13,4] = 0.2042291495073052
$arr[13,5] = 0.2147130682392819
$arr[13,6] = 3
$arr[13,7] = 1
$arr[13,8] = 60.83231733333333
$arr[13,9] = 182.496952
$arr[13,10] = 0.2679031251727568
$arr[13,11] = 0.2785759485989269
$arr[13,12] = 1638.883643774285
$arr[13,13] = 14749.95279396857
$arr[13,14] = 0.05471362740438126
$arr[13,15] = 0.05981389666134407
$arr[14,0] = 3
$arr[14,1] = 1
$arr[14,2] = 26.941003
$arr[14,3] = 80.823009
$arr[14,4] = 0.2042291495073052
$arr[14,5] = 0.2147130682392819
$arr[14,6] = 2
$arr[14,7] = 1
$arr[14,8] = 26.0983795
$arr[14,9] = 52.196759
$arr[14,10] = 0.1149362335102661
$arr[14,11] = 0.07967673702415903
$arr[14,12] = 703.1165204046384
$arr[14,13] = 4218.699122427831
$arr[14,14] = 0.02347332921737469
$arr[14,15] = 0.01710763667375157
$arr[15,0] = 3
$arr[15,1] = 1
$arr[15,2] = 33.430692
$arr[15,3] = 100.292076
$arr[15,4] = 0.2534249298216801
$arr[15,5] = 0.2664342694547198
$arr[15,6] = 3
$arr[15,7] = 1
$arr[15,8] = 10.68421466666667
$arr[15,9] = 32.052644
$arr[15,10] = 0.04705285980693976
$arr[15,11] = 0.04892736897547583
$arr[15,12] = 357.1806897832161
$arr[15,13] = 3214.626208048944
$arr[15,14] = 0.01192436769448306
$arr[15,15] = 0.01303592780932243
$arr[16,0] = 3
$arr[16,1] = 1
$arr[16,2] = 33.430692
$arr[16,3] = 100.292076
$arr[16,4] = 0.2534249298216801
$arr[16,5] = 0.2664342694547198
$arr[16,6] = 3
$arr[16,7] = 1
$arr[16,8] = 70.36235166666667
$arr[16,9] = 211.087055
$arr[16,10] = 0.3098730203341347
$arr[16,11] = 0.3222178559101571
$arr[16,12] = 2352.26210696402
$arr[16,13] = 21170.35896267618
$arr[16,14] = 0.07852954843181013
$arr[16,15] = 0.08584987904468887
$arr[17,0] = 3
$arr[17,1] = 1
$arr[17,2] = 33.430692
$arr[17,3] = 100.292076
$arr[17,4] = 0.2534249298216801
$arr[17,5] = 0.2664342694547198
$arr[17,6] = 3
$arr[17,7] = 1
$arr[17,8] = 59.09107466666666
$arr[17,9] = 177.273224
$arr[17,10] = 0.2602347611759026
$arr[17,11] = 0.2706020894912812
$arr[17,12] = 1975.455517130336
$arr[17,13] = 17779.09965417302
$arr[17,14] = 0.06594997608816479
$arr[17,15] = 0.07209767002653023
$arr[18,0] = 3
$arr[18,1] = 1
$arr[18,2] = 33.430692
$arr[18,3] = 100.292076
$arr[18,4] = 0.2534249298216801
$arr[18,5] = 0.2664342694547198
$arr[18,6] = 3
$arr[18,7] = 1
$arr[18,8] = 60.83231733333333
$arr[18,9] = 182.496952
$arr[18,10] = 0.2679031251727568
$arr[18,11] = 0.2785759485989269
$arr[18,12] = 2033.666464416928
$arr[18,13] = 18302.99817975235
$arr[18,14] = 0.06789333069591467
$arr[18,15] = 0.07422217935261066
$arr[19,0] = 3
$arr[19,1] = 1
$arr[19,2] = 33.430692
$arr[19,3] = 100.292076
$arr[19,4] = 0.2534249298216801
$arr[19,5] = 0.2664342694547198
$arr[19,6] = 2
$arr[19,7] = 1
$arr[19,8] = 26.0983795
$arr[19,9] = 52.196759
$arr[19,10] = 0.1149362335102661
$arr[19,11] = 0.07967673702415903
$arr[19,12] = 872.4868867636141
$arr[19,13] = 5234.921320581684
$arr[19,14] = 0.02912770691130743
$arr[19,15] = 0.02122861322156764
$arr[20,0] = 2
$arr[20,1] = 1
$arr[20,2] = 19.3233515
$arr[20,3] = 38.646703
$arr[20,4] = 0.1464827290385481
$arr[20,5] = 0.1026681916589156
$arr[20,6] = 3
$arr[20,7] = 1
$arr[20,8] = 10.68421466666667
$arr[20,9] = 32.052644
$arr[20,10] = 0.04705285980693976
$arr[20,11] = 0.04892736897547583
$arr[20,12] = 206.4548355054553
$arr[20,13] = 1238.729013032732
$arr[20,14] = 0.006892431313588747
$arr[20,15] = 0.005023284495340633
$arr[21,0] = 2
$arr[21,1] = 1
$arr[21,2] = 19.3233515
$arr[21,3] = 38.646703
$arr[21,4] = 0.1464827290385481
$arr[21,5] = 0.1026681916589156
$arr[21,6] = 3
$arr[21,7] = 1
$arr[21,8] = 70.36235166666667
$arr[21,9] = 211.087055
$arr[21,10] = 0.3098730203341347
$arr[21,11] = 0.3222178559101571
$arr[21,12] = 1359.636453621611
$arr[21,13] = 8157.818721729665
$arr[21,14] = 0.04539104567396156
$arr[21,15] = 0.03308152458650886
$arr[22,0] = 2
$arr[22,1] = 1
$arr[22,2] = 19.3233515
$arr[22,3] = 38.646703
$arr[22,4] = 0.1464827290385481
$arr[22,5] = 0.1026681916589156
$arr[22,6] = 3
$arr[22,7] = 1
$arr[22,8] = 59.09107466666666
$arr[22,9] = 177.273224
$arr[22,10] = 0.2602347611759026
$arr[22,11] = 0.2706020894912812
$arr[22,12] = 1141.837606296745
$arr[22,13] = 6851.025637780473
$arr[22,14] = 0.03811989800774102
$arr[22,15] = 0.02778222718719389
$arr[23,0] = 2
$arr[23,1] = 1
$arr[23,2] = 19.3233515
$arr[23,3] = 38.646703
$arr[23,4] = 0.1464827290385481
$arr[23,5] = 0.1026681916589156
$arr[23,6] = 3
$arr[23,7] = 1
$arr[23,8] = 60.83231733333333
$arr[23,9] = 182.496952
$arr[23,10] = 0.2679031251727568
$arr[23,11] = 0.2785759485989269
$arr[23,12] = 1175.484250391543
$arr[23,13] = 7052.905502349256
$arr[23,14] = 0.03924318089326117
$arr[23,15] = 0.02860088888231884
$arr[24,0] = 2
$arr[24,1] = 1
$arr[24,2] = 19.3233515
$arr[24,3] = 38.646703
$arr[24,4] = 0.1464827290385481
$arr[24,5] = 0.1026681916589156
$arr[24,6] = 2
$arr[24,7] = 1
$arr[24,8] = 26.0983795
$arr[24,9] = 52.196759
$arr[24,10] = 0.1149362335102661
$arr[24,11] = 0.07967673702415903
$arr[24,12] = 504.3081606588943
$arr[24,13] = 2017.232642635577
$arr[24,14] = 0.01683617314999561
$arr[24,15] = 0.008180266507553376
$ws.Range("E2:T26").Value = $arr
